$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header J1: 입찰공고번호 -> 비고
$ws.Range("J1").Value = "비고"

# Add new rows 121-127
$data = @(
    @{ Row=121; A="당진중학교 외 2교 전기차충전시설 설치 전기공사"; B="충청남도당진교육지원청"; C="2026-01-16"; D=29592000; E=29525874; F=26656255; G=87.745; H=90.28100000000001; I=99.77654095701541 }
    @{ Row=122; A="구미교육지원청 청사 남측 부출입구 신설 전기공사"; B="경상북도구미교육지원청"; C="2026-01-16"; D=59631000; E=59678642; F=53747182; G=87.745; H=90.06100000000001; I=100.0798946856501 }
    @{ Row=123; A="여좌천 일원 보행등 보수 전기공사"; B="경상남도 창원시 진해구"; C="2026-01-15"; D=90264000; E=85919373; F=77525910; G=87.745; H=90.23099999999999; I=95.18675551714969 }
    @{ Row=124; A="군서초 후관동 화장실개조 및 급식실동 창고 증축 전기공사"; B="전라남도영광교육지원청"; C="2026-01-16"; D=35207000; E=35176184; F=31736657; G=87.745; H=90.22199999999999; I=99.91247195160054 }
    @{ Row=125; A="순천부영초 외 4교 내진보강 외 3건 전기공사"; B="전라남도순천교육지원청"; C="2026-01-16"; D=48978000; E=49241546; F=44459700; G=87.745; H=90.289; I=100.5380905712769 }
    @{ Row=126; A="해양수련원 조리실 환기설비개선 전기공사"; B="충청남도보령교육지원청"; C="2026-01-15"; D=33114000; E=33070595; F=29869031; G=87.745; H=90.319; I=99.86892251011656 }
    @{ Row=127; A="광주자연과학고 노후 급식실 환경개선 전기공사 감리용역"; B="광주광역시교육청"; C="2026-01-15"; D=14760428; E=14760428; F=13290290; G=87.745; H=90.04000000000001; I=100 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B

    # Force column C (date-like text) to remain a text string, not auto-convert to a date serial
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $item.C
    $cCell.Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = $item.I

    # Column J: empty placeholder cell (mirrors empty J cells used throughout the sheet)
    $ws.Range("J120").Copy($ws.Cells.Item($r, 10))
}
